$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "managingOffice"
$ws.Range("B20").Value = "Glasgow"

$ws.Range("B20").Select()
